$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format A7 as text first so the date-like string "2023-10-09" is stored
# as a literal string (matching the other date cells in column A) instead
# of being auto-converted into a date serial number.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2023-10-09"
# Drop the temporary text number-format again so the cell ends up back on
# the workbook's default style (same as every other cell in the sheet).
$ws.Range("A7").ClearFormats()

$ws.Range("B7").Value = "This is a title"
$ws.Range("C7").Value = "Entry "
